$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zipcode_CheckOrder")

$ws.Range("E1").Value = "Feedback (Optional)"
$ws.Range("E2").Value = "This hidden cell has failed for checking Columbus zipcode"
$ws.Range("E4").Value = "This hidden cell has failed for checking Mountain View zipcode"
$ws.Range("E3").Value = "This hidden cell has failed for checking New York zipcode"

$ws.Columns.Item(5).ColumnWidth = 54.6666666666667

$ws.Range("E3").Select()
